$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Ricordeau" reference text: a stray mid-paragraph line break between
# "...on a simple" and "example which permitted..." is joined into one line.
# Re-assigning this (slightly corrected) text causes the shared-string table
# to drop the old entry and append the corrected one at the end - matching
# the reorder of <si> elements seen in the workbook XML.
$fixedText = @'
used FCA to generalize policies in reinforcement learning by grouping similar states using their descriptions.
%\\2003 paper
\\One of the very interesting properties of Reinforcement Learning algorithms is that they allow learning without prior knowledge of the environment. However, when the agents use algorithms that enable a generalization of the learning, they are unable to explain their choices. Neural networks are good examples of this problem. After a reminder about the basis of Reinforcement Learning, the Lattice Concept will be introduced. Then,
Q-Concept-Learning, a Reinforcement Learning algorithm that enables a generalization of the learning, the use of structured languages as well as an explanation of the agent’s choices will be presented.
\begin{comment}
\\First, we exposed the basis of Reinforcement Learning and concept lattice. We explained how to use concept lattice to structure the memory of an agent which
learns with reinforcement and the algorithmic consequences.
\\The first advantage of this method is the possibility of using more structured languages than attributes-value languages to describe the environment of the agent, using the Learning-Languages (Liqui`ere M. (2002) [4]).
\\Then, we presented Q-Concept-Learning, a Reinforcement Learning algorithm which allows generalizing the learning, i.e. to use previously learned knowledge on new situations. This algorithm keeps the symbolic representation, contrary to other generalization
methods like approximation function or neural networks. Moreover, Q-Concept-Learning allows building a part of the model of the agent’s environment. It can return the useful concepts for the task it learned. We implemented Q-Concept-Learning on a simple example which permitted to show the improvement due to generalization as well as the return of the useful concepts.
%\\2007
\\The generalization of policies in reinforcement learning is a
main issue, both from the theoretical model point of view and for their
applicability. However, generalizing from a set of examples or searching
for regularities is a problem which has already been intensively studied
in machine learning. Our work uses techniques in which generalizations
are constrained by a language bias, in order to regroup similar states.
Such generalizations are principally based on the properties of concept
lattices. To guide the possible groupings of similar environment’s states,
we propose a general algebraic framework, considering the generalization
of policies through a set partition of the states and using a language bias
as an a priori knowledge. We give an application as an example of our
approach by proposing and experimenting a bottom-up algorithm.
\\First, we reminded the reader of the general principles in the use of a Galois lattice as a generalization space for a set of objects described with a generalization
language. Then, we extended these results to define an original generalization
space for the partitions instead of powerset, in particular by introducing the
notion of partition agreeing with a language. We showed that, considering generalization as a partitioning, this structure can be used to generalize learning
in reinforcement learning. Finally, we proposed a bottom-up algorithm which
through a reinforcement learning and a generalization language made it possible
to produce general rules linking description elements of the environment and
optimal actions.
This work considered principally the algebraical aspects. Consequently, we didn’t
treat two remaining important issues: “How to use the rules on unknown examples” and “How to deal with knowledge uncertainty”.
Our future work will focus on three directions. First, a more formal link between
our work and the relational reinforcement learning [9]. Secondly, an algorithmic improvement, studying the different ways of exploring the space search and
the classical balancing between spatial and temporal complexity. Thirdly, we
will consider the exploration-exploitation problem, i.e. using previously acquired
knowledge or acquiring new knowledge and as a corollary, managing knowledge
uncertainty. It remains a difficult problem in machine learning and in reinforcement learning in particular. In our case, it’s transposed in a particular way as
the concept reliability, already presented in [10]. We think that a relevant approach can consists in the explicit management of second order elements (states
distribution, maximal reward,. . . ) directly through the algorithms.
\end{comment}
'@
$ws.Range("B14").Value2 = $fixedText

# The row had grown taller to fit the (now single-line-joined) paragraph.
$ws.Rows.Item(1).RowHeight = 375

# Reset the view back to the top-left cell (no special scrolled/selected state).
$ws.Range("A1").Select()
